# Update the Week 16 distribution parameters logged on the YDS sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(1.0636703576051187, 0.8198896503734165, 0.695308674575605, 3.287717118945786)"
$ws.Range("C2").Value = "JSU(-1.6693635800731657, 1.3310910126843543, 1.1454092494056438, 4.700573310041975)"
$ws.Range("D2").Value = "JSU(-0.785754844881186, 0.9472243134278244, 1.2660135120402045, 2.077152262094618)"
$ws.Range("E2").Value = "NIG(1.7532180036021794, 1.2282916145317602, 4.2803936983605215, 6.854969068638511)"
